# Unify the conception of DataNode, DataTable, Entity.
# This workbook's sheet was renamed from "Property1" to "DataNode" to match
# the unified naming scheme, and the view's active selection moved to H13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab/name.
$ws.Name = "DataNode"

# Update the stored selection (active cell) on the sheet, matching the
# frozen-pane "bottomLeft" pane that was active in the source file.
$ws.Range("H13").Select()
